# Update the ObjTables/SBtab header metadata (date + objTablesVersion)
# embedded in the descriptive header strings that live in cell A1 (and,
# for the very first sheet, also A2) of every worksheet.
#
#   date='2020-04-27 01:08:41'  -> date='2020-05-29 00:22:36'
#   date='2020-04-27 01:08:42'  -> date='2020-05-29 00:22:36'
#   objTablesVersion='0.0.9'    -> objTablesVersion='1.0.0'

$wb = $excel.ActiveWorkbook

function Update-ObjTablesCell($range) {
    $text = $range.Text
    if ($text -eq $null -or $text -eq "") {
        return
    }
    if ($text -like "*ObjTables*") {
        $updated = $text -replace "date='2020-04-27 \d{2}:\d{2}:\d{2}'", "date='2020-05-29 00:22:36'"
        $updated = $updated -replace "objTablesVersion='0\.0\.9'", "objTablesVersion='1.0.0'"
        if ($updated -ne $text) {
            $range.Value = $updated
        }
    }
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $wasProtected = $ws.ProtectContents
    if ($wasProtected) {
        $ws.Unprotect()
    }

    Update-ObjTablesCell $ws.Range("A1")
    Update-ObjTablesCell $ws.Range("A2")

    if ($wasProtected) {
        $ws.Protect()
    }
}
